$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the added date column (BY), matching the text style used
# by the existing date headers in row 1 (column BX etc.)
$ws.Range("BY1").Value = "26-oct"

# New counts for each product row, matching the numeric/centered style
# used by the existing count columns (e.g. BX2:BX11)
$values = @{
    2  = 6
    3  = 8
    4  = 7
    5  = 8
    6  = 5
    7  = 6
    8  = 14
    9  = 6
    10 = 10
    11 = 1
}

foreach ($row in $values.Keys) {
    $cell = $ws.Range("BY$row")
    $cell.Value = $values[$row]
    $cell.HorizontalAlignment = -4108
    $cell.NumberFormat = "0"
}

# Reproduce the final selection left behind in the workbook (one row below
# the last entered value, as if the user had just finished typing the column)
$ws.Range("BY12").Select()
